$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 164, shifting the existing rows
# 164-261 down to 165-262 (dimension grows from A1:T261 to A1:T262).
$ws.Rows("164:164").Insert()

# Populate the newly inserted row 164 with the new weekly record.
$ws.Cells.Item(164, 1).Value = 7
$ws.Cells.Item(164, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(164, 3).Value = "Ñuble"
$ws.Cells.Item(164, 4).Value = 44879
$ws.Cells.Item(164, 5).Value = 16
$ws.Cells.Item(164, 6).Value = "Fruta"
$ws.Cells.Item(164, 7).Value = 100108
$ws.Cells.Item(164, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(164, 9).Value = 100108005
$ws.Cells.Item(164, 10).Value = "Piña"
$ws.Cells.Item(164, 11).Value = "Caramelo"
$ws.Cells.Item(164, 12).Value = "Segunda"
$ws.Cells.Item(164, 13).Value = 40
$ws.Cells.Item(164, 14).Value = 25000
$ws.Cells.Item(164, 15).Value = 25000
$ws.Cells.Item(164, 16).Value = 25000
$ws.Cells.Item(164, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(164, 18).Value = "Ecuador"
$ws.Cells.Item(164, 19).Value = 1786
$ws.Cells.Item(164, 20).Value = 14
